$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" '29.582.86'
$ws.Range("E2").Value = '  -2.64%  '

# Row 3
Set-TextCell "D3" '2.000.88'
$ws.Range("E3").Value = '  -5.08%  '

# Row 4
Set-TextCell "D4" '1.015'
$ws.Range("E4").Value = '  +0.82%  '

# Row 5
Set-TextCell "D5" '331.07'
$ws.Range("E5").Value = '  -3.96%  '

# Row 6
$ws.Range("E6").Value = '  +0.69%  '

# Row 7
Set-TextCell "D7" '0.5000'
$ws.Range("E7").Value = '  -4.44%  '

# Row 8
Set-TextCell "D8" '0.4246'
$ws.Range("E8").Value = '  -4.49%  '

# Row 9
Set-TextCell "D9" '54.46'
$ws.Range("E9").Value = '  -0.63%  '

# Row 10
Set-TextCell "D10" '0.09142'
$ws.Range("E10").Value = '  -2.62%  '

# Row 11
Set-TextCell "D11" '1.122'
$ws.Range("E11").Value = '  -4.41%  '

# Row 12
Set-TextCell "D12" '23.43'
$ws.Range("E12").Value = '  -6.15%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell "D13" '2.077.67'
$ws.Range("E13").Value = '  +0.88%  '

# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell "D14" '8.091'
$ws.Range("E14").Value = '  -6.60%  '

# Row 15
Set-TextCell "D15" '6.522'
$ws.Range("E15").Value = '  -6.02%  '

# Row 16
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell "D16" '1.015'
$ws.Range("E16").Value = '  +0.69%  '

# Row 17
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell "D17" '94.95'
$ws.Range("E17").Value = '  -6.85%  '

# Row 18
Set-TextCell "D18" '0.00001117'
$ws.Range("E18").Value = '  -3.90%  '

# Row 19
Set-TextCell "D19" '0.06675'
$ws.Range("E19").Value = '  -0.64%  '

# Row 20
Set-TextCell "D20" '19.85'
$ws.Range("E20").Value = '  -6.50%  '

# Row 21
Set-TextCell "D21" '1.011'
$ws.Range("E21").Value = '  +0.64%  '

# Row 22
Set-TextCell "D22" '5.963'
$ws.Range("E22").Value = '  -6.13%  '

# Row 23
Set-TextCell "D23" '29.619.58'
$ws.Range("E23").Value = '  -2.63%  '

# Row 24
Set-TextCell "D24" '12.04'
$ws.Range("E24").Value = '  -4.78%  '

# Row 25
Set-TextCell "D25" '2.282'
$ws.Range("E25").Value = '  -0.65%  '

# Row 26
Set-TextCell "D26" '159.06'
$ws.Range("E26").Value = '  -2.36%  '

# Row 27
Set-TextCell "D27" '20.75'
$ws.Range("E27").Value = '  -5.78%  '

# Row 28
Set-TextCell "D28" '6.390'
$ws.Range("E28").Value = '  -6.03%  '

# Row 29
Set-TextCell "D29" '2.314'
$ws.Range("E29").Value = '  -8.43%  '

# Row 30
Set-TextCell "D30" '128.86'
$ws.Range("E30").Value = '  -3.73%  '

# Row 31
Set-TextCell "D31" '1.062'
$ws.Range("E31").Value = '  -7.56%  '

# Row 32
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell "D32" '0.09936'
$ws.Range("E32").Value = '  -5.73%  '

# Row 33
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell "D33" '1.572'
$ws.Range("E33").Value = '  -9.91%  '

# Row 34
Set-TextCell "D34" '5.853'
$ws.Range("E34").Value = '  -6.64%  '

# Row 35
Set-TextCell "D35" '3.813'
$ws.Range("E35").Value = '  -2.83%  '

# Row 36
Set-TextCell "D36" '9.525'
$ws.Range("E36").Value = '  -8.02%  '

# Row 37
Set-TextCell "D37" '0.02475'
$ws.Range("E37").Value = '  -6.21%  '

# Row 38
Set-TextCell "D38" '1.318'
$ws.Range("E38").Value = '  -3.28%  '

# Row 39
Set-TextCell "D39" '0.06382'
$ws.Range("E39").Value = '  -6.16%  '

# Row 40
Set-TextCell "D40" '0.6595'
$ws.Range("E40").Value = '  -6.61%  '

# Row 41
Set-TextCell "D41" '11.74'
$ws.Range("E41").Value = '  -6.76%  '

# Row 42
Set-TextCell "D42" '0.2066'
$ws.Range("E42").Value = '  -7.30%  '

# Row 43
Set-TextCell "D43" '1.012'
$ws.Range("E43").Value = '  +0.69%  '

# Row 44
Set-TextCell "D44" '0.6362'
$ws.Range("E44").Value = '  -7.47%  '

# Row 45
Set-TextCell "D45" '13.56'
$ws.Range("E45").Value = '  -5.95%  '

# Row 46
Set-TextCell "D46" '2.212'
$ws.Range("E46").Value = '  -6.40%  '

# Row 47
Set-TextCell "D47" '1.287'
$ws.Range("E47").Value = '  -7.54%  '

# Row 48
Set-TextCell "D48" '3.530'
$ws.Range("E48").Value = '  -3.21%  '

# Row 49
Set-TextCell "D49" '0.00000000338'
$ws.Range("E49").Value = '  -5.30%  '

# Row 50
Set-TextCell "D50" '0.06989'
$ws.Range("E50").Value = '  -3.59%  '

# Row 51
Set-TextCell "D51" '1.126'
$ws.Range("E51").Value = '  -7.07%  '
